$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J, matching the style used by the
# other header cells in row 1 (bold, centered, bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I ("I0") and J ("IF"), rows 2-22.
$iValues = @(6, 8, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
$jValues = @(6, 8, 1, 3, 6, 6, 5, 4, 6, 6, 7, 6, 7, 6, 5, 6, 5, 3, 6, 4, 2)

for ($idx = 0; $idx -lt 21; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
